$excel = New-Object -ComObject Excel.Application
$wb2 = $excel.ActiveWorkbook
$ws = $wb2.ActiveSheet

# The edit rotates the data of rows 12, 13 and 14 (the row number / A:AY
# cell formatting stays put, only the field values move):
#   new row12 <- old row13
#   new row13 <- old row14
#   new row14 <- old row12
# Cell-by-cell Range.Copy is used (rather than Value assignment) so that
# text that merely looks like a date ("2026-02-16") or time ("10:33")
# round-trips as plain text instead of being re-interpreted by Excel as a
# date/time serial, and so that cells which are simply absent in the
# source row stay absent in the destination row instead of turning into
# empty placeholder cells.

$row12Cols = @("A","B","D","E","F","G","H","I","K","L","M","N","P","Q","R","S","T","U","V","W","Y","AA","AC","AD","AE","AG","AT","AW","AX","AY")
$row13Cols = @("A","B","D","E","F","G","H","I","J","K","N","P","Q","R","S","T","U","V","W","Y","Z","AA","AB","AD","AE","AF","AG","AT","AW","AX","AY")
$row14Cols = @("A","B","D","E","F","G","H","I","J","K","N","P","Q","R","S","T","U","V","W","Y","AA","AC","AD","AE","AF","AG","AT","AW","AX","AY")

$fullRange = "A:AY"
$scratchRow = 999

function Copy-RowCells($cols, $srcRow, $dstRow) {
    foreach ($c in $cols) {
        $srcRef = "$c$srcRow"
        $dstRef = "$c$dstRow"
        $ws.Range($srcRef).Copy($ws.Range($dstRef)) | Out-Null
    }
}

function Clear-Row($r) {
    $ws.Range("A" + $r + ":AY" + $r).ClearContents() | Out-Null
}

# 1) stash old row 12 in a scratch row
Copy-RowCells $row12Cols 12 $scratchRow

# 2) row 12 <- old row 13
Clear-Row 12
Copy-RowCells $row13Cols 13 12

# 3) row 13 <- old row 14
Clear-Row 13
Copy-RowCells $row14Cols 14 13

# 4) row 14 <- old row 12 (from scratch)
Clear-Row 14
Copy-RowCells $row12Cols $scratchRow 14

# 5) wipe the scratch row so it doesn't linger in the sheet
Clear-Row $scratchRow

# 6) row 21: observer list order changed
$ws.Range("AX21").Value = "Lars-Erik Nilsson, Anna-Lena Thommson"
